$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet's Hyperlinks.Delete() only works when called on the whole
# sheet's collection (it is not properly scoped when called through a
# single-cell Range), so the reliable way to drop just one hyperlink
# (the old mailto: link that lived on E2) is to remove every hyperlink and
# recreate the ones that should remain, in their original order, so the
# relationship ids line up the same way as before. A brand new hyperlink for
# J2 (youtube) is added last, reusing the same mailto target the old E2
# hyperlink used to point to.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("M2"), "https://lightingstores.com.sa/en")
$ws.Range("M2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("O2"), "https://lightingstores.com.sa/en")
$ws.Range("O2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("P2"), "https://i.imgur.com/otJ9G8X.png")
$ws.Range("P2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("L2"), "http://www.lightingstores.com/")
$ws.Range("L2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("N2"), "http://www.iluslighting.com/")
$ws.Range("N2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("R2"), "https://rmoosa2014.github.io/Resume/LS LOGO.png,https://rmoosa2014.github.io/Resume/Illuslogo.svg,https://rmoosa2014.github.io/Resume/Illictlogo.svg,https://rmoosa2014.github.io/Resume/HYP_Logo.png,")
$ws.Range("R2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C2"), "https://rmoosa2014.github.io/Resume/LS LOGO.png")
$ws.Range("C2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("G2"), "https://sa.linkedin.com/in/riyas-moosa-72923992?trk=people-guest_people_search-card")
$ws.Range("G2").Style = "Hyperlink"

# emailAddress (E2): new email, no longer a hyperlink
$ws.Range("E2").Value = "ahmed.maher@example.com"
$ws.Range("E2").Style = "Normal"

# youtube (J2): new value, now carries the (reused) mailto hyperlink
$ws.Range("J2").Value = "https://youtube.com/@ahmed"
$ws.Hyperlinks.Add($ws.Range("J2"), "mailto:Riyasmoosa@example.com")
$ws.Range("J2").Style = "Hyperlink"

# Update the active selection
$ws.Range("I5").Select()
